$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a text value to a cell while avoiding Excel's automatic
# numeric/locale parsing (important for comma-separated coordinate strings
# and decimal-looking confidence values), and without leaving a stray
# cell style applied afterwards.
function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# Row 16
$ws.Range("D16").Value = "image_20250807111314_ppp0.jpg"
Set-TextValue $ws.Range("I16") "643,531,686,575"
Set-TextValue $ws.Range("J16") "0.76"

# Row 17
$ws.Range("D17").Value = "image_20250807111314_ppp0.jpg"
Set-TextValue $ws.Range("I17") "794,481,830,526"
Set-TextValue $ws.Range("J17") "0.72"

# Row 18
$ws.Range("D18").Value = "image_20250808100711_ppp0.jpg"
Set-TextValue $ws.Range("I18") "1182,409,1232,451"
Set-TextValue $ws.Range("J18") "0.75"
